$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create a new "Author" worksheet (seed data for a CMS author) and insert it
# as the very first sheet, before the existing "Blogs" sheet.
# ---------------------------------------------------------------------------
$wsAuthor = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$wsAuthor.Name = "Author"

$wsAuthor.Range("A1").Value = "Name"
$wsAuthor.Range("B1").Value = "Email"
$wsAuthor.Range("B2").Value = "cms_author@domain.com"

# ---------------------------------------------------------------------------
# Add an "AuthorName" column to "BlogPosts" (right after "BlogName") and
# stamp every existing blog post row with the new seed author.
# ---------------------------------------------------------------------------
$wsBlogPosts = $wb.Worksheets.Item("BlogPosts")
$wsBlogPosts.Columns("B:B").Insert()
$wsBlogPosts.Range("B1").Value = "AuthorName"
$wsBlogPosts.Range("B2").Value = "cms_author"
$wsBlogPosts.Range("B3").Value = "cms_author"
$wsBlogPosts.Range("B4").Value = "cms_author"
$wsBlogPosts.Range("B5").Value = "cms_author"

# Fill in the author's "Name" (re-uses the shared "cms_author" string above).
$wsAuthor.Range("A2").Value = "cms_author"

$wsAuthor.Columns("A:B").AutoFit()
$wsBlogPosts.Columns("B:B").AutoFit()

# ---------------------------------------------------------------------------
# Restore selection / active-sheet state: Author sheet leaves A2 selected,
# BlogPosts stays the active tab with C11 selected.
# ---------------------------------------------------------------------------
[void]$wsAuthor.Range("A2").Select()
[void]$wsBlogPosts.Activate()
[void]$wsBlogPosts.Range("C11").Select()
